$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "54.385.37"
$ws.Range("E2").Value = "  +0.85%  "
Set-TextValue "D3" "2.285.05"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.28%  "
Set-TextValue "D5" "503.98"
$ws.Range("E5").Value = "  +1.88%  "
Set-TextValue "D6" "129.67"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  +0.08%  "
Set-TextValue "D9" "0.0955"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("E10").Value = "  +1.28%  "
Set-TextValue "D11" "0.334"
$ws.Range("E11").Value = "  +3.37%  "
Set-TextValue "D12" "4.73"
$ws.Range("E12").Value = "  +1.49%  "
Set-TextValue "D13" "2.691.65"
$ws.Range("E13").Value = "  +0.51%  "
Set-TextValue "D14" "23.06"
$ws.Range("E14").Value = "  +6.97%  "
Set-TextValue "D15" "54.337.40"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("E16").Value = "  +0.43%  "
Set-TextValue "D17" "2.274.14"
$ws.Range("E17").Value = "  -0.83%  "
Set-TextValue "D18" "10.28"
$ws.Range("E18").Value = "  +3.90%  "
$ws.Range("E19").Value = "  +2.37%  "
Set-TextValue "D20" "305.43"
$ws.Range("E20").Value = "  +2.14%  "
Set-TextValue "D21" "6.41"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("E22").Value = "  +0.43%  "
Set-TextValue "D23" "62.11"
$ws.Range("E23").Value = "  -2.25%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  +2.26%  "
$ws.Range("E26").Value = "  +2.90%  "
Set-TextValue "D27" "174.07"
$ws.Range("E27").Value = "  +6.90%  "
$ws.Range("E28").Value = "  +1.10%  "
Set-TextValue "D29" "6.00"
$ws.Range("E29").Value = "  +3.19%  "
Set-TextValue "D30" "0.0₃0692"
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("E32").Value = "  +0.01%  "
Set-TextValue "D33" "17.81"
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("E34").Value = "  +0.07%  "
Set-TextValue "D35" "0.939"
$ws.Range("E35").Value = "  +9.38%  "
$ws.Range("E36").Value = "  +1.87%  "
$ws.Range("E37").Value = "  +3.71%  "
$ws.Range("E38").Value = "  -0.53%  "
Set-TextValue "D39" "1.42"
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("E40").Value = "  +1.97%  "
Set-TextValue "D41" "5.02"
$ws.Range("E41").Value = "  +2.06%  "
Set-TextValue "D42" "125.04"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("E43").Value = "  +3.55%  "
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("E46").Value = "  +1.13%  "
Set-TextValue "D47" "0.373"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("E49").Value = "  +1.02%  "
Set-TextValue "D50" "16.43"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("E51").Value = "  +0.21%  "
